$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142, shifting the existing rows 142:266 down to 143:267.
$ws.Rows("142:142").Insert()

# Populate the newly-inserted row 142 with this week's price quote.
$ws.Range("A142").Value = 3
$ws.Range("B142").Value = "Femacal de La Calera"
$ws.Range("C142").Value = "Coquimbo"
$ws.Range("D142").Value = 44566
$ws.Range("E142").Value = 5
$ws.Range("F142").Value = 100114013
$ws.Range("G142").Value = "Zanahoria"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 310
$ws.Range("K142").Value = 6500
$ws.Range("L142").Value = 7000
$ws.Range("M142").Value = 6742
$ws.Range("N142").Value = "`$/saco 20 kilos"
$ws.Range("O142").Value = "Provincia de Quillota"
$ws.Range("P142").Value = 337
$ws.Range("Q142").Value = 20
$ws.Range("R142").Value = "Hortaliza"
